$wb = $excel.ActiveWorkbook

# --- Rename worksheets ---
$wsTimeSeries = $wb.Worksheets.Item(1)
$wsTimeSeries.Name = "TimeSeries"

$wsSummary = $wb.Worksheets.Item(2)
$wsSummary.Name = "Summary"

$wsFountain = $wb.Worksheets.Item(3)
$wsFountain.Name = "FountainManeuvers&ChasingEps"

# --- Fix capitalization of a label on the Summary sheet ---
$wsSummary.Range("B7").Value = "The euclidean distance between dolphin and fish"

# --- Widen column B on the Summary sheet (target stored width ~38.332) ---
$wsSummary.Columns.Item(2).ColumnWidth = 37.5

# --- Update section header labels on the FountainManeuvers&ChasingEps sheet ---
$wsFountain.Range("A2").Value = "Fountain maneuvers before learning"
$wsFountain.Range("A9").Value = "Fountain maneuvers after learning"
$wsFountain.Range("A14").Value = "Chasing episodes before learning"

# --- Update selections on each sheet ---
$wsTimeSeries.Activate()
$wsTimeSeries.Range("I22").Select()

$wsSummary.Activate()
$wsSummary.Range("H17").Select()

$wsFountain.Activate()
$wsFountain.Range("A14").Select()
